$d = $word.ActiveDocument

# Locate the final paragraph in the document (currently holds a single
# space character) - this is the paragraph whose run gets removed.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastRange = $lastPara.Range

# Insert all of the new paragraphs right before the final paragraph mark,
# using a collapsed range positioned just before it so Word appends new
# sibling paragraphs after the existing (still-intact) last paragraph
# rather than merging into it.
$insertionPoint = $d.Range($lastRange.End - 1, $lastRange.End - 1)

$newParagraphsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="hu-HU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="hu-HU"/></w:rPr><w:t xml:space="preserve">Hasonlóan járunk el az adatbázisba töltés esetén is. Azért, hogy demonstrálni tudjuk a feladatkiírásban jelölt műveleteket, nem ad </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="hu-HU"/></w:rPr><w:t>hozzá</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="hu-HU"/></w:rPr><w:t xml:space="preserve"> hogy 14 merőpanel adatait mutassuk, elég ha 2 szerepel. Amennyiben az összesre szükség van, mindössze hozzá kell fűzni a többi táblát a lekérdezésbe, mivel az adatbázisban szerepelnek. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="160" w:line="278" w:lineRule="auto"/><w:rPr><w:lang w:val="hu-HU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="hu-HU"/></w:rPr><w:br w:type="page"/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="hu-HU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="hu-HU"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Az adatok betöltése megtörtént az </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="hu-HU"/></w:rPr><w:t>adabázisba</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="hu-HU"/></w:rPr><w:t xml:space="preserve">. Ettől a ponttól kilépünk a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="hu-HU"/></w:rPr><w:t>python</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="hu-HU"/></w:rPr><w:t xml:space="preserve"> világából és megérkezünk a tiszta SQL műveletek világába.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="hu-HU"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="hu-HU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="hu-HU"/></w:rPr><w:t xml:space="preserve">Az adatbázis műveletek sorban </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="hu-HU"/></w:rPr><w:t>ismertetem(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="hu-HU"/></w:rPr><w:t>ahol szükséges) és a megállapításokat SQL commentként írom.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="hu-HU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="hu-HU"/></w:rPr><w:t>A csatolt adatbázis fájlban ellenőrizhető az elkészített táblák és annak tartalmai.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="hu-HU"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="hu-HU"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="hu-HU"/></w:rPr></w:pPr></w:p>
'@

$insertionPoint.InsertXML($newParagraphsXml)

# Now that the new paragraphs exist after it, clear the original last
# paragraph's run so it becomes an empty paragraph (just the mark).
$clearRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$clearRange.Text = ""
